# Update survey results:
#  1. Fix some wording inconsistencies in existing responses.
#  2. Add a newly received survey response ("Miss Strong").
#  3. Re-sort all responses alphabetically by respondent name.
#  4. Leave the sheet selection on the cell the editor ended up on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Wording / typo corrections on existing rows -------------------------
$ws.Range("E8").Value  = "Programmer"
$ws.Range("E9").Value  = "Film editor"
$ws.Range("F11").Value = "It is helpful for the future of technology "
$ws.Range("F14").Value = "Helps you develop better problem solving skills"
$ws.Range("E16").Value = "Computer engineer"

# --- 2. Append the new survey response ---------------------------------------
$ws.Range("A23").Value = "Miss Strong"
$ws.Range("B23").Value = 8
$ws.Range("C23").Value = 40
$ws.Range("D23").Value = "Not trolling"
$ws.Range("E23").Value = "Software developer"
$ws.Range("F23").Value = "You learn how to think critically and systematically"

# Give the new row the same look (style/height) as the rest of the data rows.
$ws.Range("A22:F22").Copy()
$ws.Range("A23:F23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Rows.Item(23).RowHeight = 15.75

# --- 3. Sort all responses (A2:F23) alphabetically by name (column A) -------
$sortRange = $ws.Range("A2:F23")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A23"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws.Sort.Apply()

# --- 4. Restore the selected cell -------------------------------------------
$ws.Range("F3").Select()
